$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet ---
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("C2").Value  = 6.165366141587581
$ws.Range("C3").Value  = 20.187170153299732
$ws.Range("C4").Value  = -0.7561407398332505
$ws.Range("C6").Value  = 8.316041303807811
$ws.Range("C7").Value  = 26.802829663958995
$ws.Range("C8").Value  = 0.448111843819463
$ws.Range("C10").Value = 8.316041303807811
$ws.Range("C11").Value = 26.802829663958995
$ws.Range("C12").Value = 0.448111843819463
$ws.Range("C14").Value = 7.277951412170038
$ws.Range("C15").Value = 23.609576965842123
$ws.Range("C16").Value = 0.32270298420759413
$ws.Range("C18").Value = 7.206623920809382
$ws.Range("C19").Value = 23.39016754275452
$ws.Range("C20").Value = -0.0164820496491173

# --- LANDING GEARS sheet ---
$ws2 = $wb.Worksheets.Item("LANDING GEARS")
$ws2.Range("C2").Value = 18.01244712946852
